$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.57465587734637
$ws.Range("C2").Value = 9.256536877775636
$ws.Range("D2").Value = 7.852815787582225
$ws.Range("E2").Value = 9.589500852679581
$ws.Range("F2").Value = 41.57193039689169
$ws.Range("L2").Value = 10.48871128386846
$ws.Range("N2").Value = 20.10905685934477
$ws.Range("B3").Value = 24.07473069002785
$ws.Range("C3").Value = 8.639776147237725
$ws.Range("D3").Value = 7.875895845758193
$ws.Range("E3").Value = 9.608231250050478
$ws.Range("F3").Value = 41.20086126237158
$ws.Range("L3").Value = 10.47213174224126
$ws.Range("N3").Value = 20.16594502128074
$ws.Range("B4").Value = 23.77074027902023
$ws.Range("C4").Value = 8.237740905599779
$ws.Range("D4").Value = 7.891315933408241
$ws.Range("E4").Value = 9.620485003207421
$ws.Range("F4").Value = 40.98430932995671
$ws.Range("L4").Value = 10.46427795857787
$ws.Range("N4").Value = 20.20286001205958
$ws.Range("B5").Value = 23.64779423517631
$ws.Range("C5").Value = 8.068018890434338
$ws.Range("D5").Value = 7.897912594987202
$ws.Range("E5").Value = 9.625668289068519
$ws.Range("F5").Value = 40.89897459810187
$ws.Range("L5").Value = 10.46166356297643
$ws.Range("N5").Value = 20.21840202707703
$ws.Range("B6").Value = 23.62744093717553
$ws.Range("C6").Value = 8.039480391074255
$ws.Range("D6").Value = 7.899026826828052
$ws.Range("E6").Value = 9.626540443350624
$ws.Range("F6").Value = 40.88498273707314
$ws.Range("L6").Value = 10.46126486766312
$ws.Range("N6").Value = 20.22101288555701
$ws.Range("B7").Value = 23.76907816234872
$ws.Range("C7").Value = 8.235475847662437
$ws.Range("D7").Value = 7.891403633108585
$ws.Range("E7").Value = 9.620554137878532
$ws.Range("F7").Value = 40.9831465926424
$ws.Range("L7").Value = 10.46424032567529
$ws.Range("N7").Value = 20.20306759773537
$ws.Range("B8").Value = 24.40177773233025
$ws.Range("C8").Value = 9.048710301975357
$ws.Range("D8").Value = 7.860513881002144
$ws.Range("E8").Value = 9.595802992537996
$ws.Range("F8").Value = 41.44168969227987
$ws.Range("L8").Value = 10.48251270127107
$ws.Range("N8").Value = 20.12825957886546
$ws.Range("B9").Value = 25.65793799425832
$ws.Range("C9").Value = 10.45915852073315
$ws.Range("D9").Value = 7.809901798322059
$ws.Range("E9").Value = 9.553225437205827
$ws.Range("F9").Value = 42.426891027911
$ws.Range("L9").Value = 10.53673197986323
$ws.Range("N9").Value = 19.99733227498857
$ws.Range("B10").Value = 26.57971266673561
$ws.Range("C10").Value = 11.38431239429836
$ws.Range("D10").Value = 7.778865465389805
$ws.Range("E10").Value = 9.525553289585883
$ws.Range("F10").Value = 43.1980181505291
$ws.Range("L10").Value = 10.58765641627452
$ws.Range("N10").Value = 19.91077829836113
$ws.Range("B11").Value = 26.99670009313971
$ws.Range("C11").Value = 11.78133141734204
$ws.Range("D11").Value = 7.766097839834954
$ws.Range("E11").Value = 9.513743412342867
$ws.Range("F11").Value = 43.55791197821225
$ws.Range("L11").Value = 10.61319629520936
$ws.Range("N11").Value = 19.87350149933703
$ws.Range("B12").Value = 27.15410227271434
$ws.Range("C12").Value = 11.92826517657324
$ws.Range("D12").Value = 7.761458749946483
$ws.Range("E12").Value = 9.509382875959517
$ws.Range("F12").Value = 43.69539951388216
$ws.Range("L12").Value = 10.62320523640169
$ws.Range("N12").Value = 19.85968796943779
$ws.Range("B13").Value = 27.12022777162949
$ws.Range("C13").Value = 11.89677166229859
$ws.Range("D13").Value = 7.762449132073526
$ws.Range("E13").Value = 9.510317036637773
$ws.Range("F13").Value = 43.66573728793914
$ws.Range("L13").Value = 10.62103468546538
$ws.Range("N13").Value = 19.86264949774936
$ws.Range("B14").Value = 27.00966064485426
$ws.Range("C14").Value = 11.79348795262514
$ws.Range("D14").Value = 7.765712247859116
$ws.Range("E14").Value = 9.513382433378762
$ws.Range("F14").Value = 43.56919965180963
$ws.Range("L14").Value = 10.61401299556393
$ws.Range("N14").Value = 19.87235898672656
$ws.Range("B15").Value = 26.94186494043002
$ws.Range("C15").Value = 11.72978034865761
$ws.Range("D15").Value = 7.767736534610759
$ws.Range("E15").Value = 9.51527460281817
$ws.Range("F15").Value = 43.51022113002936
$ws.Range("L15").Value = 10.6097558441045
$ws.Range("N15").Value = 19.87834573411007
$ws.Range("B16").Value = 26.55240017975763
$ws.Range("C16").Value = 11.35788822344565
$ws.Range("D16").Value = 7.779727169317105
$ws.Range("E16").Value = 9.52634072768754
$ws.Range("F16").Value = 43.17467246422963
$ws.Range("L16").Value = 10.58603477912572
$ws.Range("N16").Value = 19.91325669267002
$ws.Range("B17").Value = 26.31275740521247
$ws.Range("C17").Value = 11.12365259573351
$ws.Range("D17").Value = 7.787430084208512
$ws.Range("E17").Value = 9.53332855840655
$ws.Range("F17").Value = 42.9710816176113
$ws.Range("L17").Value = 10.5720881400301
$ws.Range("N17").Value = 19.93521111340251
$ws.Range("B18").Value = 26.17471107605079
$ws.Range("C18").Value = 10.98668155982668
$ws.Range("D18").Value = 7.791987666561869
$ws.Range("E18").Value = 9.537421051518177
$ws.Range("F18").Value = 42.85484497165243
$ws.Range("L18").Value = 10.56429012241621
$ws.Range("N18").Value = 19.94803604585928
$ws.Range("B19").Value = 26.12794009139922
$ws.Range("C19").Value = 10.93991937229957
$ws.Range("D19").Value = 7.793552565599311
$ws.Range("E19").Value = 9.538819294543462
$ws.Range("F19").Value = 42.81564071565736
$ws.Range("L19").Value = 10.56168838180647
$ws.Range("N19").Value = 19.95241221468435
$ws.Range("B20").Value = 26.3382907720735
$ws.Range("C20").Value = 11.14881957749088
$ws.Range("D20").Value = 7.786596934018987
$ws.Range("E20").Value = 9.53257711022802
$ws.Range("F20").Value = 42.99266560792123
$ws.Range("L20").Value = 10.57354965584806
$ws.Range("N20").Value = 19.93285359903672
$ws.Range("B21").Value = 27.042151792688
$ws.Range("C21").Value = 11.82391724134134
$ws.Range("D21").Value = 7.764748467713964
$ws.Range("E21").Value = 9.512479025516987
$ws.Range("F21").Value = 43.59752328345608
$ws.Range("L21").Value = 10.61606630812315
$ws.Range("N21").Value = 19.86949886090738
$ws.Range("B22").Value = 27.49916959382487
$ws.Range("C22").Value = 12.24527275301514
$ws.Range("D22").Value = 7.751610809712266
$ws.Range("E22").Value = 9.499994159851214
$ws.Range("F22").Value = 43.9997977071784
$ws.Range("L22").Value = 10.64581876654312
$ws.Range("N22").Value = 19.82985583129754
$ws.Range("B23").Value = 27.25557731152014
$ws.Range("C23").Value = 12.02219825100455
$ws.Range("D23").Value = 7.758517679910613
$ws.Range("E23").Value = 9.506598158576832
$ws.Range("F23").Value = 43.78449436641139
$ws.Range("L23").Value = 10.62976087263954
$ws.Range("N23").Value = 19.85085246293843
$ws.Range("B24").Value = 26.32674797297617
$ws.Range("C24").Value = 11.13744876870856
$ws.Range("D24").Value = 7.786973199294087
$ws.Range("E24").Value = 9.53291660601808
$ws.Range("F24").Value = 42.98290495356545
$ws.Range("L24").Value = 10.57288821872582
$ws.Range("N24").Value = 19.93391879890908
$ws.Range("B25").Value = 25.31757069663378
$ws.Range("C25").Value = 10.09719606646459
$ws.Range("D25").Value = 7.822519489547765
$ws.Range("E25").Value = 9.564108244325883
$ws.Range("F25").Value = 42.15167167374261
$ws.Range("L25").Value = 10.52010627326206
$ws.Range("N25").Value = 20.03106089177426
